# Update "想去人数" (F column) figures across sheets to reflect the
# latest generated snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 513
$ws.Range("F4").Value = 483
$ws.Range("F5").Value = 8883
$ws.Range("F6").Value = 19
$ws.Range("F7").Value = 11403
$ws.Range("F8").Value = 94
$ws.Range("F20").Value = 427
$ws.Range("F22").Value = 1937
$ws.Range("F23").Value = 743
$ws.Range("F24").Value = 676
$ws.Range("F25").Value = 369
$ws.Range("F26").Value = 298
$ws.Range("F29").Value = 1368
$ws.Range("F35").Value = 475
$ws.Range("F36").Value = 316
$ws.Range("F37").Value = 47
$ws.Range("F38").Value = 361
$ws.Range("F39").Value = 347
$ws.Range("F40").Value = 41
$ws.Range("F41").Value = 150
$ws.Range("F42").Value = 544
$ws.Range("F43").Value = 404
$ws.Range("F44").Value = 125
$ws.Range("F45").Value = 826
$ws.Range("F46").Value = 663
$ws.Range("F48").Value = 188
$ws.Range("F49").Value = 173

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 39
$ws.Range("F5").Value = 3
$ws.Range("F10").Value = 22
$ws.Range("F19").Value = 118
$ws.Range("F24").Value = 84
$ws.Range("F25").Value = 403

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2864

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 513
$ws.Range("F5").Value = 39
$ws.Range("F8").Value = 483
$ws.Range("F9").Value = 8883
$ws.Range("F10").Value = 19
$ws.Range("F11").Value = 11403
$ws.Range("F12").Value = 94
$ws.Range("F18").Value = 427
$ws.Range("F19").Value = 1937
$ws.Range("F20").Value = 743
$ws.Range("F21").Value = 676
$ws.Range("F22").Value = 369
$ws.Range("F23").Value = 298
$ws.Range("F28").Value = 1368
$ws.Range("F35").Value = 475
$ws.Range("F36").Value = 316
$ws.Range("F38").Value = 361
$ws.Range("F39").Value = 150
$ws.Range("F40").Value = 544
$ws.Range("F41").Value = 404
$ws.Range("F42").Value = 125
$ws.Range("F43").Value = 826
$ws.Range("F45").Value = 403
$ws.Range("F46").Value = 663
$ws.Range("F48").Value = 188
$ws.Range("F49").Value = 173
